$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Delivered the model in under 7 weeks, working in a 3-person team",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Delivered the model in under 7 weeks, working in a 4-person team",
    2)

Write-Output "Replaced: $found"
